# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets to reflect a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 2460
$ws1.Range("F8").Value  = 1826
$ws1.Range("F9").Value  = 3099
$ws1.Range("F11").Value = 4595
$ws1.Range("F12").Value = 420
$ws1.Range("F13").Value = 244
$ws1.Range("F14").Value = 143
$ws1.Range("F16").Value = 275
$ws1.Range("F17").Value = 628
$ws1.Range("F21").Value = 126
$ws1.Range("F23").Value = 4595
$ws1.Range("F25").Value = 22
$ws1.Range("F27").Value = 4513
$ws1.Range("F31").Value = 616
$ws1.Range("F33").Value = 48
$ws1.Range("F34").Value = 102
$ws1.Range("F35").Value = 709
$ws1.Range("F36").Value = 35
$ws1.Range("F37").Value = 650
$ws1.Range("F38").Value = 646

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 2460
$ws4.Range("F11").Value = 1826
$ws4.Range("F13").Value = 3099
$ws4.Range("F15").Value = 4595
$ws4.Range("F16").Value = 420
$ws4.Range("F17").Value = 244
$ws4.Range("F18").Value = 143
$ws4.Range("F20").Value = 275
$ws4.Range("F21").Value = 628
$ws4.Range("F26").Value = 126
$ws4.Range("F28").Value = 4595
$ws4.Range("F30").Value = 22
$ws4.Range("F32").Value = 4513
$ws4.Range("F33").Value = 8
$ws4.Range("F36").Value = 616
$ws4.Range("F39").Value = 48
$ws4.Range("F40").Value = 102
$ws4.Range("F41").Value = 709
$ws4.Range("F42").Value = 35
$ws4.Range("F43").Value = 650
$ws4.Range("F44").Value = 646
